$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93, shifting existing rows 93..202 down to 94..203
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly record
$ws.Cells.Item(93, 1).Value = 10
$ws.Cells.Item(93, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(93, 3).Value = "La Araucanía"
$ws.Cells.Item(93, 4).Value = 44483
$ws.Cells.Item(93, 4).NumberFormat = $ws.Cells.Item(94, 4).NumberFormat
$ws.Cells.Item(93, 5).Value = 9
$ws.Cells.Item(93, 6).Value = 100112009
$ws.Cells.Item(93, 7).Value = "Acelga"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 50
$ws.Cells.Item(93, 11).Value = 7000
$ws.Cells.Item(93, 12).Value = 8000
$ws.Cells.Item(93, 13).Value = 7600
$ws.Cells.Item(93, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(93, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(93, 16).Value = 633
$ws.Cells.Item(93, 17).Value = 12
$ws.Cells.Item(93, 18).Value = "Hortaliza"
